$wb = $excel.ActiveWorkbook

$wsMementos = $wb.Worksheets.Item("mementos")
$wsMissing  = $wb.Worksheets.Item("missing")
$wsSputnick = $wb.Worksheets.Item("sputnick")

# --- sheet "mementos" (sheet1): B16 becomes a hyperlink (claemit url, same as its own text) ---
$wsMementos.Hyperlinks.Add($wsMementos.Range("B16"), $wsMementos.Range("B16").Value) | Out-Null
$wsMementos.Range("B16").Style = "Hyperlink"
$wsMementos.Range("E36").Select() | Out-Null

# --- sheet "missing" (sheet2): rework the VIDEO row, note two newly-uploaded .mov files ---
$wsMissing.Range("B10").Value = "f-claemit.mov missing"
$wsMissing.Range("B11").Value = "f-room.mov missing"
$wsMissing.Range("B13").Select() | Out-Null

# --- sheet "sputnick" (sheet3): A4 becomes a hyperlink too ---
$wsSputnick.Hyperlinks.Add($wsSputnick.Range("A4"), $wsSputnick.Range("A4").Value) | Out-Null
$wsSputnick.Range("A4").Style = "Hyperlink"
$wsSputnick.Range("A15").Select() | Out-Null

# --- make "missing" the active/selected sheet (was "sputnick") ---
$wsMissing.Activate()
